$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New / renamed method labels ---
$ws.Range("B4").Value = "Holden"
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("B11").Value = "Matthies Hex"

# --- Append two new rows (28, 29) at the bottom ---
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"

# --- Rewrite the simulated data block C4:T31 (rerun of the simulation) ---
$data = New-Object "object[,]" 28,18
$data[0,0] = 0.9846718918593863
$data[0,1] = 1.070414079790856
$data[0,2] = 1.019669966659893
$data[0,3] = 0.9112654441258361
$data[0,4] = 0.9112654441258361
$data[0,5] = 1.096187508036424
$data[0,6] = 1.096187508036424
$data[0,7] = 1.01828576209981
$data[0,8] = 0.9112654441258361
$data[0,9] = 1.01828576209981
$data[0,10] = 1.057236635068117
$data[0,11] = 1.057236635068117
$data[0,12] = 1.044714412265375
$data[0,13] = 1.00857957142069
$data[0,14] = 1.00857957142069
$data[0,15] = 0.9842510395969766
$data[0,16] = 0.9842510395969766
$data[0,17] = 1.016749108762034
$data[1,0] = 1.15967593165236
$data[1,1] = 1.174750890627671
$data[1,2] = 0.6422291735971748
$data[1,3] = 1.41602622811095
$data[1,4] = 1.41602622811095
$data[1,5] = 0.3655711095286229
$data[1,6] = 0.3655711095286229
$data[1,7] = 1.065859246414155
$data[1,8] = 1.41602622811095
$data[1,9] = 1.065859246414155
$data[1,10] = 0.7157151779713889
$data[1,11] = 0.7157151779713889
$data[1,12] = 0.6912198431799842
$data[1,13] = 0.9491521946845761
$data[1,14] = 0.9491521946845761
$data[1,15] = 1.06587070304117
$data[1,16] = 1.06587070304117
$data[1,17] = 0.9706854299884892
$data[2,0] = 1.041711041305126
$data[2,1] = 1.211709231932758
$data[2,2] = 0.8888890447121203
$data[2,3] = 1.010401068960471
$data[2,4] = 1.010401068960471
$data[2,5] = 0.8511985933908125
$data[2,6] = 0.8511985933908125
$data[2,7] = 1.057405428091932
$data[2,8] = 1.010401068960471
$data[2,9] = 1.057405428091932
$data[2,10] = 0.9543020107413724
$data[2,11] = 0.9543020107413724
$data[2,12] = 0.9324976887316216
$data[2,13] = 0.9730016968144053
$data[2,14] = 0.9730016968144053
$data[2,15] = 0.9823515398509217
$data[2,16] = 0.9823515398509217
$data[2,17] = 1.010219068065537
$data[3,0] = 0.9847955213112388
$data[3,1] = 1.07733514477666
$data[3,2] = 1.015052264697405
$data[3,3] = 0.9172441617363132
$data[3,4] = 0.9172441617363132
$data[3,5] = 1.098549936276654
$data[3,6] = 1.098549936276654
$data[3,7] = 1.018573960266572
$data[3,8] = 0.9172441617363132
$data[3,9] = 1.018573960266572
$data[3,10] = 1.058561948271613
$data[3,11] = 1.058561948271613
$data[3,12] = 1.044058720413544
$data[3,13] = 1.011456019426513
$data[3,14] = 1.011456019426513
$data[3,15] = 0.9879030550039634
$data[3,16] = 0.9879030550039634
$data[3,17] = 1.018591831510807
$data[4,0] = 0.9291369369976259
$data[4,1] = 1.418330462172698
$data[4,2] = 1.057137663504084
$data[4,3] = 0.58272788125859
$data[4,4] = 0.58272788125859
$data[4,5] = 1.478359657136803
$data[4,6] = 1.478359657136803
$data[4,7] = 1.108357553965204
$data[4,8] = 0.58272788125859
$data[4,9] = 1.108357553965204
$data[4,10] = 1.293358605551003
$data[4,11] = 1.293358605551003
$data[4,12] = 1.214618291535364
$data[4,13] = 1.056481697453532
$data[4,14] = 1.056481697453532
$data[4,15] = 0.9380432434047965
$data[4,16] = 0.9380432434047965
$data[4,17] = 1.095675025839167
$data[5,0] = 1.000076061940366
$data[5,1] = 0.9857662250145853
$data[5,2] = 1.000183334274813
$data[5,3] = 1.013895530320609
$data[5,4] = 1.013895530320609
$data[5,5] = 0.9956122433398245
$data[5,6] = 0.9956122433398245
$data[5,7] = 0.9950192942077858
$data[5,8] = 1.013895530320609
$data[5,9] = 0.9950192942077858
$data[5,10] = 0.9953157687738051
$data[5,11] = 0.9953157687738051
$data[5,12] = 0.9969382906074742
$data[5,13] = 1.00150902262274
$data[5,14] = 1.00150902262274
$data[5,15] = 1.004605649547207
$data[5,16] = 1.004605649547207
$data[5,17] = 0.9984254481829972
$data[6,0] = 0.9998899103356096
$data[6,1] = 0.9968942296842329
$data[6,2] = 0.9999635757726629
$data[6,3] = 1.004022653101704
$data[6,4] = 1.004022653101704
$data[6,5] = 0.9997593175139168
$data[6,6] = 0.9997593175139168
$data[6,7] = 0.9982742615212041
$data[6,8] = 1.004022653101704
$data[6,9] = 0.9982742615212041
$data[6,10] = 0.9990167895175605
$data[6,11] = 0.9990167895175605
$data[6,12] = 0.9993323849359279
$data[6,13] = 1.000685410712275
$data[6,14] = 1.000685410712275
$data[6,15] = 1.001519721309632
$data[6,16] = 1.001519721309632
$data[6,17] = 0.9998006579882217
$data[7,0] = 0.9986395087887119
$data[7,1] = 0.9761204295953141
$data[7,2] = 1.002724914057606
$data[7,3] = 1.020327310461459
$data[7,4] = 1.020327310461459
$data[7,5] = 0.997589033456908
$data[7,6] = 0.997589033456908
$data[7,7] = 0.9916335107254579
$data[7,8] = 1.020327310461459
$data[7,9] = 0.9916335107254579
$data[7,10] = 0.9946112720911829
$data[7,11] = 0.9946112720911829
$data[7,12] = 0.9973158194133239
$data[7,13] = 1.003183284881275
$data[7,14] = 1.003183284881275
$data[7,15] = 1.007469291276321
$data[7,16] = 1.007469291276321
$data[7,17] = 0.9978391178475761
$data[8,0] = 0.9272951428562918
$data[8,1] = 1.433197042878555
$data[8,2] = 1.057913229309711
$data[8,3] = 0.5707094790691642
$data[8,4] = 0.5707094790691642
$data[8,5] = 1.492205085168194
$data[8,6] = 1.492205085168194
$data[8,7] = 1.111374110853279
$data[8,8] = 0.5707094790691642
$data[8,9] = 1.111374110853279
$data[8,10] = 1.301789598010736
$data[8,11] = 1.301789598010736
$data[8,12] = 1.220497475110395
$data[8,13] = 1.058096225030212
$data[8,14] = 1.058096225030212
$data[8,15] = 0.9362495385399502
$data[8,16] = 0.9362495385399502
$data[8,17] = 1.098782348355866
$data[9,0] = 1.038434475307687
$data[9,1] = 1.095080174008727
$data[9,2] = 0.9155725405961338
$data[9,3] = 1.044033744311262
$data[9,4] = 1.044033744311262
$data[9,5] = 0.8434491542108119
$data[9,6] = 0.8434491542108119
$data[9,7] = 1.037809204995599
$data[9,8] = 1.044033744311262
$data[9,9] = 1.037809204995599
$data[9,10] = 0.9406291796032054
$data[9,11] = 0.9406291796032054
$data[9,12] = 0.9322769666008481
$data[9,13] = 0.9750973678392242
$data[9,14] = 0.9750973678392242
$data[9,15] = 0.9923314619572334
$data[9,16] = 0.9923314619572334
$data[9,17] = 0.9957298822383701
$data[10,0] = 0.9797735730315781
$data[10,1] = 1.134844881621056
$data[10,2] = 1.016612528294735
$data[10,3] = 0.868279403494737
$data[10,4] = 0.868279403494737
$data[10,5] = 1.139126178242105
$data[10,6] = 1.139126178242105
$data[10,7] = 1.033684357694739
$data[10,8] = 0.868279403494737
$data[10,9] = 1.033684357694739
$data[10,10] = 1.086405267968422
$data[10,11] = 1.086405267968422
$data[10,12] = 1.063141021410526
$data[10,13] = 1.013696646477194
$data[10,14] = 1.013696646477194
$data[10,15] = 0.9773423357315796
$data[10,16] = 0.9773423357315796
$data[10,17] = 1.028720153729825
$data[11,0] = 1.005967438500199
$data[11,1] = 1.008690389885227
$data[11,2] = 0.9848299741513878
$data[11,3] = 1.022711986353048
$data[11,4] = 1.022711986353048
$data[11,5] = 0.9761227336360272
$data[11,6] = 0.9761227336360272
$data[11,7] = 0.9980795551939257
$data[11,8] = 1.022711986353048
$data[11,9] = 0.9980795551939257
$data[11,10] = 0.9871011444149764
$data[11,11] = 0.9871011444149764
$data[11,12] = 0.9863440876604469
$data[11,13] = 0.9989714250610003
$data[11,14] = 0.9989714250610003
$data[11,15] = 1.004906565384012
$data[11,16] = 1.004906565384012
$data[11,17] = 0.9994003462866359
$data[12,0] = 0.8725643200000008
$data[12,1] = 1.770108699999998
$data[12,2] = 1.097047600000001
$data[12,3] = 0.24511332
$data[12,4] = 0.24511332
$data[12,5] = 1.8669278
$data[12,6] = 1.8669278
$data[12,7] = 1.200521499999999
$data[12,8] = 0.24511332
$data[12,9] = 1.200521499999999
$data[12,10] = 1.533724649999999
$data[12,11] = 1.533724649999999
$data[12,12] = 1.388165633333333
$data[12,13] = 1.10418754
$data[12,14] = 1.10418754
$data[12,15] = 0.8894189849999997
$data[12,16] = 0.8894189849999997
$data[12,17] = 1.17538054
$data[13,0] = 1.134
$data[13,1] = 0.21380688
$data[13,2] = 0.39134565
$data[13,3] = 3.4428622
$data[13,4] = 3.4428622
$data[13,5] = 0.023300756
$data[13,6] = 0.023300756
$data[13,7] = 0.42160412
$data[13,8] = 3.4428622
$data[13,9] = 0.42160412
$data[13,10] = 0.222452438
$data[13,11] = 0.222452438
$data[13,12] = 0.2787501753333333
$data[13,13] = 1.295922358666667
$data[13,14] = 1.295922358666667
$data[13,15] = 1.832657319
$data[13,16] = 1.832657319
$data[13,17] = 0.9378199343333332
$data[14,0] = 1.2491109
$data[14,1] = 2.7689596
$data[14,2] = 0.20651447
$data[14,3] = 1.0822396
$data[14,4] = 1.0822396
$data[14,5] = 0.032036222
$data[14,6] = 0.032036222
$data[14,7] = 1.3724463
$data[14,8] = 1.0822396
$data[14,9] = 1.3724463
$data[14,10] = 0.7022412609999999
$data[14,11] = 0.7022412609999999
$data[14,12] = 0.5369989973333333
$data[14,13] = 0.8289073739999999
$data[14,14] = 0.828907374
$data[14,15] = 0.8922404305
$data[14,16] = 0.8922404305
$data[14,17] = 1.118551182
$data[15,0] = 1.0720652
$data[15,1] = 1.8223291
$data[15,2] = 0.66043898
$data[15,3] = 1.0505035
$data[15,4] = 1.0505035
$data[15,5] = 0.8782308599999999
$data[15,6] = 0.8782308599999999
$data[15,7] = 1.1172741
$data[15,8] = 1.0505035
$data[15,9] = 1.1172741
$data[15,10] = 0.9977524799999999
$data[15,11] = 0.9977524799999999
$data[15,12] = 0.8853146466666666
$data[15,13] = 1.015336153333333
$data[15,14] = 1.015336153333333
$data[15,15] = 1.02412799
$data[15,16] = 1.02412799
$data[15,17] = 1.10014029
$data[16,0] = 0.9883621123287667
$data[16,1] = 0.6863503221917808
$data[16,2] = 1.024913314657534
$data[16,3] = 1.284299527671233
$data[16,4] = 1.284299527671233
$data[16,5] = 0.9465206117260273
$data[16,6] = 0.9465206117260273
$data[16,7] = 0.8885736301369862
$data[16,8] = 1.284299527671233
$data[16,9] = 0.8885736301369862
$data[16,10] = 0.9175471209315067
$data[16,11] = 0.9175471209315067
$data[16,12] = 0.9533358521735159
$data[16,13] = 1.039797923178082
$data[16,14] = 1.039797923178082
$data[16,15] = 1.10092332430137
$data[16,16] = 1.10092332430137
$data[16,17] = 0.9698365864520548
$data[17,0] = 1.005530841052632
$data[17,1] = 1.853414268421052
$data[17,2] = 0.8018328494736843
$data[17,3] = 0.7597024536842105
$data[17,4] = 0.7597024536842105
$data[17,5] = 1.210274136421053
$data[17,6] = 1.210274136421053
$data[17,7] = 1.14886451
$data[17,8] = 0.7597024536842105
$data[17,9] = 1.14886451
$data[17,10] = 1.179569323210526
$data[17,11] = 1.179569323210526
$data[17,12] = 1.053657165298246
$data[17,13] = 1.039613700035088
$data[17,14] = 1.039613700035088
$data[17,15] = 0.9696358884473684
$data[17,16] = 0.9696358884473684
$data[17,17] = 1.129936509842105
$data[18,0] = 1.138605148421053
$data[18,1] = 1.192750551052632
$data[18,2] = 0.6701385515789475
$data[18,3] = 1.378430999473684
$data[18,4] = 1.378430999473684
$data[18,5] = 0.4639429711578947
$data[18,6] = 0.4639429711578947
$data[18,7] = 1.060011593157895
$data[18,8] = 1.378430999473684
$data[18,9] = 1.060011593157895
$data[18,10] = 0.7619772821578947
$data[18,11] = 0.7619772821578947
$data[18,12] = 0.7313643719649123
$data[18,13] = 0.9674618545964911
$data[18,14] = 0.9674618545964911
$data[18,15] = 1.070204140815789
$data[18,16] = 1.070204140815789
$data[18,17] = 0.9839799691403509
$data[19,0] = 0.9279295441835341
$data[19,1] = 1.166244549060898
$data[19,2] = 1.101895647005625
$data[19,3] = 0.7594890053520396
$data[19,4] = 0.7594890053520396
$data[19,5] = 1.440559545107162
$data[19,6] = 1.440559545107162
$data[19,7] = 0.992844387141834
$data[19,8] = 0.7594890053520396
$data[19,9] = 0.992844387141834
$data[19,10] = 1.216701966124498
$data[19,11] = 1.216701966124498
$data[19,12] = 1.178433193084874
$data[19,13] = 1.064297645867012
$data[19,14] = 1.064297645867012
$data[19,15] = 0.9880954857382689
$data[19,16] = 0.9880954857382689
$data[19,17] = 1.064827112975182
$data[20,0] = 1.003399005809783
$data[20,1] = 0.8368033757414323
$data[20,2] = 1.040837308069024
$data[20,3] = 1.018211072495288
$data[20,4] = 1.018211072495288
$data[20,5] = 0.8597309207966093
$data[20,6] = 0.8597309207966093
$data[20,7] = 0.9981067985513478
$data[20,8] = 1.018211072495288
$data[20,9] = 0.9981067985513478
$data[20,10] = 0.9289188596739786
$data[20,11] = 0.9289188596739786
$data[20,12] = 0.9662250091389938
$data[20,13] = 0.958682930614415
$data[20,14] = 0.9586829306144149
$data[20,15] = 0.9735649660846331
$data[20,16] = 0.9735649660846331
$data[20,17] = 0.9595147469105806
$data[21,0] = 1.010562656342278
$data[21,1] = 0.6076480114475263
$data[21,2] = 1.042499214400601
$data[21,3] = 1.226152087430778
$data[21,4] = 1.226152087430778
$data[21,5] = 0.8135648193026975
$data[21,6] = 0.8135648193026975
$data[21,7] = 0.9195313905494017
$data[21,8] = 1.226152087430778
$data[21,9] = 0.9195313905494017
$data[21,10] = 0.8665481049260496
$data[21,11] = 0.8665481049260496
$data[21,12] = 0.9251984747509002
$data[21,13] = 0.9864160990942925
$data[21,14] = 0.9864160990942925
$data[21,15] = 1.046350096178414
$data[21,16] = 1.046350096178414
$data[21,17] = 0.9366596965788805
$data[22,0] = 0.9769692938969803
$data[22,1] = 1.146773738040237
$data[22,2] = 1.001305353144287
$data[22,3] = 0.9259020730545918
$data[22,4] = 0.9259020730545918
$data[22,5] = 1.19622858816455
$data[22,6] = 1.19622858816455
$data[22,7] = 0.9980379874609863
$data[22,8] = 0.9259020730545918
$data[22,9] = 0.9980379874609863
$data[22,10] = 1.097133287812768
$data[22,11] = 1.097133287812768
$data[22,12] = 1.065190642923275
$data[22,13] = 1.04005621622671
$data[22,14] = 1.04005621622671
$data[22,15] = 1.01151768043368
$data[22,16] = 1.01151768043368
$data[22,17] = 1.040869505626939
$data[23,0] = 1.137013518901693
$data[23,1] = 0.876206697063997
$data[23,2] = 0.7303100739457405
$data[23,3] = 1.538037954437502
$data[23,4] = 1.538037954437502
$data[23,5] = 0.3456270970246729
$data[23,6] = 0.3456270970246729
$data[23,7] = 0.9671687472778617
$data[23,8] = 1.538037954437502
$data[23,9] = 0.9671687472778617
$data[23,10] = 0.6563979221512672
$data[23,11] = 0.6563979221512672
$data[23,12] = 0.6810353060827583
$data[23,13] = 0.9502779329133456
$data[23,14] = 0.9502779329133456
$data[23,15] = 1.097217938294385
$data[23,16] = 1.097217938294385
$data[23,17] = 0.9323940147752446
$data[24,0] = 1.022788678436913
$data[24,1] = 1.252766798361946
$data[24,2] = 0.9488432388393743
$data[24,3] = 0.8327873379431808
$data[24,4] = 0.8327873379431808
$data[24,5] = 0.9184771135653895
$data[24,6] = 0.9184771135653895
$data[24,7] = 1.103828430616287
$data[24,8] = 0.8327873379431808
$data[24,9] = 1.103828430616287
$data[24,10] = 1.011152772090838
$data[24,11] = 1.011152772090838
$data[24,12] = 0.9903829276736836
$data[24,13] = 0.9516976273749526
$data[24,14] = 0.9516976273749526
$data[24,15] = 0.9219700550170096
$data[24,16] = 0.9219700550170096
$data[24,17] = 1.013248599627182
$data[25,0] = 0.9925017200491179
$data[25,1] = 0.9857577545732567
$data[25,2] = 0.9886782788794876
$data[25,3] = 1.064187494965505
$data[25,4] = 1.064187494965505
$data[25,5] = 1.111936335554146
$data[25,6] = 1.111936335554146
$data[25,7] = 0.9618770025335127
$data[25,8] = 1.064187494965505
$data[25,9] = 0.9618770025335127
$data[25,10] = 1.036906669043829
$data[25,11] = 1.036906669043829
$data[25,12] = 1.020830538989049
$data[25,13] = 1.046000277684388
$data[25,14] = 1.046000277684388
$data[25,15] = 1.050547082004667
$data[25,16] = 1.050547082004667
$data[25,17] = 1.017489764425838
$data[26,0] = 1.030448369895988
$data[26,1] = 1.04109870623943
$data[26,2] = 0.9450875845562311
$data[26,3] = 1.006254783630723
$data[26,4] = 1.006254783630723
$data[26,5] = 0.908474684501443
$data[26,6] = 0.908474684501443
$data[26,7] = 1.041072458202626
$data[26,8] = 1.006254783630723
$data[26,9] = 1.041072458202626
$data[26,10] = 0.9747735713520345
$data[26,11] = 0.9747735713520345
$data[26,12] = 0.9648782424201
$data[26,13] = 0.9852673087782641
$data[26,14] = 0.9852673087782641
$data[26,15] = 0.9905141774913788
$data[26,16] = 0.9905141774913788
$data[26,17] = 0.99540609783774
$data[27,0] = 1.018597256835603
$data[27,1] = 1.05515999306867
$data[27,2] = 0.9103740932803251
$data[27,3] = 1.185448799620363
$data[27,4] = 1.185448799620363
$data[27,5] = 0.9887668686133244
$data[27,6] = 0.9887668686133244
$data[27,7] = 0.9550628296041742
$data[27,8] = 1.185448799620363
$data[27,9] = 0.9550628296041742
$data[27,10] = 0.9719148491087493
$data[27,11] = 0.9719148491087493
$data[27,12] = 0.9514012638326079
$data[27,13] = 1.043092832612621
$data[27,14] = 1.043092832612621
$data[27,15] = 1.078681824364556
$data[27,16] = 1.078681824364556
$data[27,17] = 1.01890164017041
$ws.Range("C4:T31").Value = $data
